$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry over the same formatting used by the rest of the table (column A)
# into the new column G, so the new cells share style index 1.
$ws.Range("A1:A6").Copy()
$ws.Range("G1:G6").PasteSpecial(-4122)

# Add a new "City" column header in G1
$ws.Range("G1").Value = "City"

# Fill the City column for each data row with the same value as the
# existing "Location" column (E)
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 5).Value2
}
